# Apply the SCD0026-012 renaming edit to the "View agenda pada Portal" workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet (and book title) from SCD0338 to SCD0026
$ws.Name = "SCD0026"

# Update the test case id cells
$ws.Range("B2").Value = "SCD0026-012"
$ws.Range("B3").Value = "SCD0026-012"

# Update the comment text cell
$ws.Range("S3").Value = "Test Komentar 4"

# Update the view: scroll position + active selection
$ws.Application.ActiveWindow.ScrollColumn = 8
$ws.Range("S4").Select()
